$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2020" column (Q), mirroring the formatting
# of the preceding "2019" column (P) for both the header row and the data row.
$ws.Range("P4:P5").Copy($ws.Range("Q4:Q5"))
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 3.3

# Move/restore the active selection to match the new state of the sheet.
$ws.Range("R4").Select()
